$wb = $excel.ActiveWorkbook
$targetSheets = @(1, 4)

foreach ($sheetIdx in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    # Row 2
    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = "2024.01.21"
    $ws.Range("C2").Value = "苏州·世纪幻想动漫游戏展"
    $ws.Range("D2").Value = "清禾路886号 尹山湖大剧院"
    $ws.Range("E2").Value = "2024.01.21 10:30-01.21 17:00"
    $ws.Range("F2").Value = 1930
    $ws.Range("G2").NumberFormat = "@"
    $ws.Range("G2").Value = "60"
    $ws.Range("H2").Value = $true
    $ws.Range("I2").Value = "https://show.bilibili.com/platform/detail.html?id=80053&msource=Msearch_colligation"
    $ws.Range("J2").Value = "//i1.hdslb.com/bfs/openplatform/202312/vtGcfnyc1703060683812.jpeg"

    # Row 3
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2024.01.27"
    $ws.Range("C3").Value = "昆山·“不是！你有病吧！”主题展（取消）"
    $ws.Range("D3").Value = "绿地大道258号游站未来城2号楼 魔之塔"
    $ws.Range("E3").Value = "2024.01.27 10:00-01.27 19:00"
    $ws.Range("F3").Value = 270
    $ws.Range("G3").NumberFormat = "@"
    $ws.Range("G3").Value = "不可售"
    $ws.Range("H3").Value = $false
    $ws.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=79124&msource=Msearch_colligation"
    $ws.Range("J3").Value = "//i2.hdslb.com/bfs/openplatform/202311/Z7mV6VXN1701160508967.jpeg"

    # Row 4
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024.01.28"
    $ws.Range("C4").Value = "苏州.第二届THO 赤维极陵"
    $ws.Range("D4").Value = "白塔东路60号(近平江路) 苏州书香府邸平江府"
    $ws.Range("E4").Value = "2024.01.28 10:00-01.28 21:00"
    $ws.Range("F4").Value = 259
    $ws.Range("G4").NumberFormat = "@"
    $ws.Range("G4").Value = "58"
    $ws.Range("H4").Value = $false
    $ws.Range("I4").Value = "https://show.bilibili.com/platform/detail.html?id=79002&msource=Msearch_colligation"
    $ws.Range("J4").Value = "//i0.hdslb.com/bfs/openplatform/202311/5AgvDWGQ1700817845950.jpeg"

    # Row 5
    $ws.Range("B5").NumberFormat = "@"
    $ws.Range("B5").Value = "2024.02.03"
    $ws.Range("C5").Value = "苏州·TCD国潮动漫游戏嘉年华"
    $ws.Range("D5").Value = "苏州大道东688号 苏州国际博览中心"
    $ws.Range("E5").Value = "2024.02.03 09:30-02.04 17:00"
    $ws.Range("F5").Value = 8754
    $ws.Range("G5").NumberFormat = "@"
    $ws.Range("G5").Value = "60"
    $ws.Range("H5").Value = $false
    $ws.Range("I5").Value = "https://show.bilibili.com/platform/detail.html?id=80084&msource=Msearch_colligation"
    $ws.Range("J5").Value = "//i0.hdslb.com/bfs/openplatform/202401/aDe3s9MS1705479547745.jpeg"

    # Row 6
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2024.02.04"
    $ws.Range("C6").Value = "苏州·TCD国潮动漫游戏嘉年华吴磊内场"
    $ws.Range("D6").Value = "苏州大道东688号 苏州国际博览中心"
    $ws.Range("E6").Value = "2024.02.04 09:30-02.04 17:00"
    $ws.Range("F6").Value = 578
    $ws.Range("G6").NumberFormat = "@"
    $ws.Range("G6").Value = "已售罄"
    $ws.Range("H6").Value = $false
    $ws.Range("I6").Value = "https://show.bilibili.com/platform/detail.html?id=80398&msource=Msearch_colligation"
    $ws.Range("J6").Value = "//i1.hdslb.com/bfs/openplatform/202401/bHsHJ3f21704186294427.jpeg"

    # Row 7
    $ws.Range("B7").NumberFormat = "@"
    $ws.Range("B7").Value = "2024.02.14"
    $ws.Range("C7").Value = "常熟·CDW·动漫展02"
    $ws.Range("D7").Value = "常熟国际展览中心 国际展览中心"
    $ws.Range("E7").Value = "2024.02.14 09:00-02.15 17:30"
    $ws.Range("F7").Value = 629
    $ws.Range("G7").NumberFormat = "@"
    $ws.Range("G7").Value = "55"
    $ws.Range("H7").Value = $false
    $ws.Range("I7").Value = "https://show.bilibili.com/platform/detail.html?id=80504&msource=Msearch_colligation"
    $ws.Range("J7").Value = "//i1.hdslb.com/bfs/openplatform/202401/VHHzVjad1704438989848.jpeg"

    # Row 8
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = "2024.02.14"
    $ws.Range("C8").Value = "常熟·漫魂动漫游戏展01"
    $ws.Range("D8").Value = "虞山北路258号 星程酒店(长江路店)"
    $ws.Range("E8").Value = "2024.02.14 09:00-02.14 21:00"
    $ws.Range("F8").Value = 87
    $ws.Range("G8").NumberFormat = "@"
    $ws.Range("G8").Value = "50"
    $ws.Range("H8").Value = $false
    $ws.Range("I8").Value = "https://show.bilibili.com/platform/detail.html?id=80248&msource=Msearch_colligation"
    $ws.Range("J8").Value = "//i2.hdslb.com/bfs/openplatform/202312/oPrKUOby1703664065719.jpeg"

    # Row 9
    $ws.Range("B9").NumberFormat = "@"
    $ws.Range("B9").Value = "2024.02.14"
    $ws.Range("C9").Value = "苏州·第一届寒假动漫展宅舞比赛-CF01"
    $ws.Range("D9").Value = "润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店"
    $ws.Range("E9").Value = "2024.02.14 10:00-02.14 16:00"
    $ws.Range("F9").Value = 21
    $ws.Range("G9").NumberFormat = "@"
    $ws.Range("G9").Value = "49"
    $ws.Range("H9").Value = $true
    $ws.Range("I9").Value = "https://show.bilibili.com/platform/detail.html?id=80528&msource=Msearch_colligation"
    $ws.Range("J9").Value = "//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg"

    # Row 10
    $ws.Range("B10").NumberFormat = "@"
    $ws.Range("B10").Value = "2024.02.16"
    $ws.Range("C10").Value = "太仓·龙狮新春动漫节4.0"
    $ws.Range("D10").Value = "滨河路126号 凯景世纪大酒店"
    $ws.Range("E10").Value = "2024.02.16 08:30-02.16 15:00"
    $ws.Range("F10").Value = 4
    $ws.Range("G10").NumberFormat = "@"
    $ws.Range("G10").Value = "45"
    $ws.Range("H10").Value = $false
    $ws.Range("I10").Value = "https://show.bilibili.com/platform/detail.html?id=81044&msource=Msearch_colligation"
    $ws.Range("J10").Value = "//i1.hdslb.com/bfs/openplatform/202401/AMDXVltp1705568031796.jpeg"

    # Row 11
    $ws.Range("B11").NumberFormat = "@"
    $ws.Range("B11").Value = "2024.02.16"
    $ws.Range("C11").Value = "苏州·Good Jump ACG迎新特别篇X动漫品牌博览会"
    $ws.Range("D11").Value = "金山南路288号 广电国际会展中心"
    $ws.Range("E11").Value = "2024.02.16 10:00-02.17 17:00"
    $ws.Range("F11").Value = 9323
    $ws.Range("G11").NumberFormat = "@"
    $ws.Range("G11").Value = "60"
    $ws.Range("H11").Value = $false
    $ws.Range("I11").Value = "https://show.bilibili.com/platform/detail.html?id=79303&msource=Msearch_colligation"
    $ws.Range("J11").Value = "//i2.hdslb.com/bfs/openplatform/202312/C3P0Encm1701659824998.jpeg"

    # Row 12
    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = "2024.02.25"
    $ws.Range("C12").Value = "苏州·第五届次元鹿角动漫游戏展"
    $ws.Range("D12").Value = "清禾路886号 尹山湖大剧院"
    $ws.Range("E12").Value = "2024.02.25 10:00-02.25 17:00"
    $ws.Range("F12").Value = 2406
    $ws.Range("G12").NumberFormat = "@"
    $ws.Range("G12").Value = "68"
    $ws.Range("H12").Value = $true
    $ws.Range("I12").Value = "https://show.bilibili.com/platform/detail.html?id=79333&msource=Msearch_colligation"
    $ws.Range("J12").Value = "//i1.hdslb.com/bfs/openplatform/202401/tqrMA6qB1704787264871.jpeg"

    # Row 13
    $ws.Range("B13").NumberFormat = "@"
    $ws.Range("B13").Value = "2024.03.08"
    $ws.Range("C13").Value = "苏州·国风宠物-cosplay展"
    $ws.Range("D13").Value = "木渎金山南路288号 苏州国际影视娱乐城"
    $ws.Range("E13").Value = "2024.03.08 09:00-03.10 17:30"
    $ws.Range("F13").Value = 16
    $ws.Range("G13").NumberFormat = "@"
    $ws.Range("G13").Value = "65"
    $ws.Range("H13").Value = $true
    $ws.Range("I13").Value = "https://show.bilibili.com/platform/detail.html?id=80635&msource=Msearch_colligation"
    $ws.Range("J13").Value = "//i2.hdslb.com/bfs/openplatform/202401/Rfd9PcBN1704781416369.jpeg"

    # Row 14
    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "2024.04.13"
    $ws.Range("C14").Value = "苏州·绘时国乙1.0-秩序之外"
    $ws.Range("D14").Value = "石路步行街永福桥浜15号 银河广场"
    $ws.Range("E14").Value = "2024.04.13 13:30-04.13 20:00"
    $ws.Range("F14").Value = 48
    $ws.Range("G14").NumberFormat = "@"
    $ws.Range("G14").Value = "78"
    $ws.Range("H14").Value = $false
    $ws.Range("I14").Value = "https://show.bilibili.com/platform/detail.html?id=80789&msource=Msearch_colligation"
    $ws.Range("J14").Value = "//i0.hdslb.com/bfs/openplatform/202401/SjKfDxBh1705041298410.jpeg"

    # Row 15
    $ws.Range("B15").NumberFormat = "@"
    $ws.Range("B15").Value = "2024.04.21"
    $ws.Range("C15").Value = "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0"
    $ws.Range("D15").Value = "清禾路888号2号楼3楼 格莱美婚礼宴会中心"
    $ws.Range("E15").Value = "2024.04.21 10:00-04.21 21:00"
    $ws.Range("F15").Value = 341
    $ws.Range("G15").NumberFormat = "@"
    $ws.Range("G15").Value = "48"
    $ws.Range("H15").Value = $true
    $ws.Range("I15").Value = "https://show.bilibili.com/platform/detail.html?id=78666&msource=Msearch_colligation"
    $ws.Range("J15").Value = "//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg"

    # Row 16
    $ws.Range("B16").NumberFormat = "@"
    $ws.Range("B16").Value = "2024.05.01"
    $ws.Range("C16").Value = "昆山·第十二届理想乡动漫游戏展"
    $ws.Range("D16").Value = "花桥经济开发区绿地大道1598号 花桥国际博览中心"
    $ws.Range("E16").Value = "2024.05.01 10:00-05.03 17:00"
    $ws.Range("F16").Value = 10394
    $ws.Range("G16").NumberFormat = "@"
    $ws.Range("G16").Value = "预售中"
    $ws.Range("H16").Value = $true
    $ws.Range("I16").Value = "https://show.bilibili.com/platform/detail.html?id=77196&msource=Msearch_colligation"
    $ws.Range("J16").Value = "//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png"

    # Row 17
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("B17").Value = "2024.05.01"
    $ws.Range("C17").Value = "苏州·第十七届 I COME ACG  动漫品牌博览会"
    $ws.Range("D17").Value = "金山南路288号 广电国际会展中心"
    $ws.Range("E17").Value = "2024.05.01 10:00-05.02 17:00"
    $ws.Range("F17").Value = 10630
    $ws.Range("G17").NumberFormat = "@"
    $ws.Range("G17").Value = "65"
    $ws.Range("H17").Value = $true
    $ws.Range("I17").Value = "https://show.bilibili.com/platform/detail.html?id=79789&msource=Msearch_colligation"
    $ws.Range("J17").Value = "//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg"

}

"done"